$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the "Prospect" territory code for THE HEARTH (row 19) ---
# Was "040", should be "023".
$ws.Range("C19").Value = "023"

# --- Add new prospect row: MSP COMMERCIAL ---
# Insert a blank row at row 24 (pushes the existing rows 24-25 down to 25-26),
# matching the formatting (row height, cell styles) of the surrounding rows.
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).RowHeight = 13.05

$ws.Range("A24").Value = "MSP COMMERCIAL"
$ws.Range("B24").Value = "Steiner, Owen A"
$ws.Range("C24").Value = "015"
# D24 (Last Invoice Date) stays blank - this is a new prospect, not yet invoiced.
$ws.Range("E24").Value = "0008352"
